$d = $word.ActiveDocument

# --- Hunk 1: remove the _GoBack bookmark that currently sits in the
#     "Renta" paragraph (after ", monto") ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete() | Out-Null
}

# --- Hunk 2: append a new "PagoTemporal" paragraph after the last
#     paragraph ("Cliente {...}.") and move the _GoBack bookmark there,
#     right before the closing "}." run ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter() | Out-Null
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="-5"/>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>PagoTemporal</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> {</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>idPago</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve">, </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>tipoPago</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve">, fecha, monto, </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>idAlumno</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>, idProfesor</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:t>}.</w:t>
  </w:r>
</w:p>
'@

$newPara.Range.InsertXML($xml) | Out-Null
